$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (column D) are written back as
# literal text, matching the source data (inline strings), instead of
# being auto-converted to numbers by the COM Value setter.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.173.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.15%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.669.69'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E4').Value = '  -0.47%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5235'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.19%  '

$ws.Range('E7').Value = '  -0.46%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2623'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.71%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06335'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.40%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.16%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07533'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.29%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.679.52'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.447'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.46%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5511'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.25%  '

$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008009'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.84%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.82%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.175.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.33%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.768'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '187.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.21%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.208'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.40%  '

$ws.Range('E23').Value = '  -0.48%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1248'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.88%  '

$ws.Range('E26').Value = '  -3.48%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.95%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06381'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.06%  '

$ws.Range('E29').Value = '  -1.76%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.276'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.02%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.520'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.40%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.418'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.06%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.649'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.90%  '

$ws.Range('E34').Value = '  -1.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6040'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.20%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.407'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.46%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.756'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('E38').Value = '  -0.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.112.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.38%  '

$ws.Range('E40').Value = '  -0.20%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8662'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.47%  '

$ws.Range('E42').Value = '  -0.70%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.35'
$ws.Range('D43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.823.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.50%  '

$ws.Range('E45').Value = '  +0.33%  '

$ws.Range('E46').Value = '  -3.02%  '

$ws.Range('E47').Value = '  -0.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.088'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.08%  '

$ws.Range('E49').Value = '  -0.75%  '

$ws.Range('E50').Value = '  -1.10%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.937'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.38%  '
